$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "Variables": rename a few variable names ---
$ws1.Range("B2").Value = "id"
$ws1.Range("B10").Value = "employ"

# Delete the "AGE_HYP" row entirely (row 24), shifting subsequent rows up
$ws1.Rows.Item(24).Delete()

# After the shift, rows 25-40 become rows 24-39.
# Rename remaining variable names that changed casing/whitespace.
$ws1.Range("B29").Value = "med_stat"
$ws1.Range("B30").Value = "med_nsaid"
$ws1.Range("B31").Value = "f1_htn_kora"
$ws1.Range("B38").Value = "f1_untdat "
$ws1.Range("B39").Value = "f2_untdat "

# --- Sheet2 "Categories": mirror the variable renames in column A ---
$ws2.Range("A31").Value = "employ"
$ws2.Range("A32").Value = "employ"
$ws2.Range("A33").Value = "employ"
$ws2.Range("A34").Value = "employ"

$ws2.Range("A69").Value = "med_stat"
$ws2.Range("A70").Value = "med_stat"

$ws2.Range("A71").Value = "med_nsaid"
$ws2.Range("A72").Value = "med_nsaid"

$ws2.Range("A73").Value = "f1_htn_kora"
$ws2.Range("A74").Value = "f1_htn_kora"

# Category label corrections
$ws2.Range("C56").Value = "Yes"
$ws2.Range("C64").Value = "I don't know"
